$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1 - copy style from an existing header cell (E1)
# so they reuse the same bold/bordered/centered header style, then set text.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# New boolean "Outliers_MAD" columns F,G,H for data rows 2-18 (default FALSE)
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}

# Flag the two detected outlier rows in the KNN_Outliers_MAD column
$ws.Cells.Item(4, 6).Value = $true
$ws.Cells.Item(9, 6).Value = $true
